$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark that currently sits on the
#    "Network: Build" heading paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Remove the "Meeting Locations: Build" heading paragraph together with
#    the paragraph describing it.
$pHeading = $d.Paragraphs.Item(4)
$pBody = $d.Paragraphs.Item(5)
$rngRemove = $d.Range($pHeading.Range.Start, $pBody.Range.End)
$rngRemove.Delete()

# 3. Rename the "Rejecting Matches: Build" heading to "Tabbed Interface".
$d.Content.Find.Execute("Rejecting Matches: Build", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Tabbed Interface", 2) | Out-Null

# 4. Replace the paragraph text that used to describe rejecting matches with
#    the new "Tabbed Interface" blurb.
$d.Content.Find.Execute("This will contain criteria for rejecting matches from the automatic webpage opener" + [char]8217 + "s list", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Look into upgrading the application to use tabs rather than a single window.", 2) | Out-Null

# 5. Re-create the "_GoBack" bookmark, now collapsed right after the text of
#    the final paragraph (before its paragraph mark). Adding a bookmark at a
#    collapsed position that sits immediately before a paragraph mark gets
#    mis-anchored by this host, so we work around it by temporarily padding
#    the paragraph with one extra character, anchoring the bookmark just
#    before that pad, and then deleting the pad again.
$pLast = $d.Paragraphs.Last
$padRange = $d.Range($pLast.Range.End - 1, $pLast.Range.End - 1)
$padRange.InsertAfter("X")

$bmRange = $d.Range($pLast.Range.End - 2, $pLast.Range.End - 2)
$d.Bookmarks.Add("_GoBack", $bmRange)

$padChar = $d.Range($d.Bookmarks.Item("_GoBack").End, $d.Bookmarks.Item("_GoBack").End + 1)
$padChar.Delete()
